$d = $word.ActiveDocument
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("paramentro", $true, $false, $false, $false, $false, $true, 1, $false, "parámetro", 2)
